$wb = $excel.ActiveWorkbook

$wsCompanies = $wb.Worksheets.Item("Companies")
$wsBios      = $wb.Worksheets.Item("Bios")
$wsNotes     = $wb.Worksheets.Item("Notes")

# ----- Companies sheet (row 2 = the single company record) -----
$wsCompanies.Range("A2").Value = 97
$wsCompanies.Range("B2").Value = "Ultimate Fight Zone Association (UFZA)"
$wsCompanies.Range("C2").Value = "UFZA"
$wsCompanies.Range("D2").Value = "www.ultimatefightzoneassociation(ufza).c"
$wsCompanies.Range("I2").Value = "ultimatefightzoneassociation(ufza)."
$wsCompanies.Range("J2").Value = "ultimatefightzoneassociation(ufza)B"
$wsCompanies.Range("K2").Value = "ultimatefightzoneassociation(u"
$wsCompanies.Range("M2").Value = 37
$wsCompanies.Range("R2").Value = 43

# ----- Bios sheet -----
$wsBios.Range("A2").Value = 97

$bio = @"
Name: Ultimate Fight Zone Association (UFZA)
Founder: John Smith
Established in: 2005
Location: Orlando, Florida
Size: Medium
About UFZA:
Ultimate Fight Zone Association (UFZA) is a professional wrestling company based in Orlando, Florida. Established in 2005 by founder John Smith, UFZA has quickly become known for its intense and action-packed wrestling matches. The company is considered medium in size and has a dedicated fan base that spans across the United States.
UFZA has a roster of talented wrestlers, including both seasoned veterans and up-and-coming stars. The company prides itself on providing top-notch entertainment for its fans, with high-energy matches and captivating storylines that keep audiences on the edge of their seats.
UFZA hosts live events in various cities across the country, drawing in crowds of enthusiastic fans who come to witness the thrilling action firsthand. The company also produces a weekly television show that showcases the best matches and moments from their live events.
In addition to live events and television programming, UFZA offers merchandise featuring their logo and the likenesses of their top wrestlers. Fans can purchase t-shirts, hats, posters, and more to show their support for their favorite UFZA stars.
Ultimate Fight Zone Association is committed to providing a platform for talented wrestlers to showcase their skills and entertain fans with their athleticism and charisma. With a focus on high-quality wrestling and engaging storytelling, UFZA continues to grow its fan base and solidify its reputation as a premier professional wrestling company.
"@
$wsBios.Range("B2").Value = $bio

# ----- Notes sheet -----
$wsNotes.Range("A2").Value = "Ultimate Fight Zone Association (UFZA)"
$wsNotes.Range("B2").Value = "a sports wrestling company"
